$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that changes from
# 45203 (2023-10-04) to 45205 (2023-10-06) for every data row (2..388).
$ws.Range("C2:C388").Value = 45205
